$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The task "Remove columns that are not needed" (old row 15) is being folded
# into task 1's ("Import and clean up metadata") Details text, so delete
# that whole row - this shifts rows 16:23 up to 15:22.
$ws.Rows("15").Delete()

# Row 14 ("Import and clean up metadata") now carries the merged Details
# text (its own original text plus the text that used to live in the
# deleted row).
$ws.Range("D14").Value = "Load the spreadsheet Global Superstore into R and perform basic model cleanup. Missing value, Nas, duplicate values etc. Refer to the Data Dictionary. Remove columns that are not needed. Only Ship Data, Market, Segment, Sales, Quantity and profit are required"
$ws.Rows("14").RowHeight = 60

# Renumber the Task Ref column (B) for the now-shifted rows 15:22 (2..9).
$ws.Range("B15").Value = 2
$ws.Range("B16").Value = 3
$ws.Range("B17").Value = 4
$ws.Range("B18").Value = 5
$ws.Range("B19").Value = 6
$ws.Range("B20").Value = 7
$ws.Range("B21").Value = 8
$ws.Range("B22").Value = 9

# Add the new "Status" column header, matching the formatting used by the
# rest of the header row (B13:E13).
$ws.Range("E13").Copy()
$ws.Range("F13").PasteSpecial(-4122)
$ws.Range("F13").Value = "Status"

# Selection / active cell moved to B23 in the saved file.
$ws.Range("B23").Select()
